# Daily attendance processing - 2026-01-09 16:40:12
# Rotate the "Recorded By" (column G) comma-separated list of recorders
# left by one position (move the first entry to the end) for every row
# that has more than one recorder listed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    $val = $cell.Value()

    if ($null -eq $val) { continue }

    $text = [string]$val
    $parts = $text -split ', '

    if ($parts.Count -gt 1) {
        $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ', '
        $cell.Value = $rotated
    }
}
